# Update column F (dSF) values for specific rows as per the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = 6
$ws.Range("F14").Value = -2
$ws.Range("F17").Value = 4
$ws.Range("F20").Value = 4
$ws.Range("F24").Value = -1
$ws.Range("F32").Value = 1
$ws.Range("F36").Value = -2
$ws.Range("F38").Value = 4
$ws.Range("F43").Value = 0
